# Kayıt silindi: 11166374
# Remove the record row whose "Kayıt No" (column A) equals 11166374 from
# both the master "Kayitlar" list and its per-department mirror sheet
# "Merkez İlçe"; all following rows shift up by one.

$wb = $excel.ActiveWorkbook

$targetId = "11166374"

$sheetNames = @("Kayitlar", "Merkez İlçe")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    for ($r = 1; $r -le $lastRow; $r++) {
        $cellText = $ws.Cells.Item($r, 1).Text
        if ($cellText -eq $targetId) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}
